$p = $ppt.ActivePresentation

# -----------------------------------------------------------------------
# 1. Table style change: the table on the "SOURCES OF FINANCE" slide gets
#    re-pointed from the (custom) default table style to a different
#    (built-in) table style GUID.
# -----------------------------------------------------------------------
$oldStyleId = "{B0F85A5B-569C-4128-8DB3-7B256CEDF2EE}"
$newStyleId = "{DFAB9F1F-DBD9-447B-A805-E55CCFD6D994}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}

# -----------------------------------------------------------------------
# 2. Re-colour theme: the presentation's design swaps from the "Integral"
#    colour scheme to the stock "Office Theme" colour scheme (fonts and
#    format scheme are already identical between the two, only the
#    colours differ).
# -----------------------------------------------------------------------
function ColorLong([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($idx = 1; $idx -le $officeThemeColors.Length; $idx++) {
    $tcs.Colors($idx).RGB = ColorLong $officeThemeColors[$idx - 1]
}
